$wb = $excel.ActiveWorkbook

# ---- Neg_Change (sheet 1) ----
$ws = $wb.Worksheets.Item(1)

# Row 2: TRENT
$ws.Range("A2").Value = "TRENT"
$ws.Range("B2").Value = 3906
$ws.Range("C2").Value = 3950
$ws.Range("D2").Value = 3886.2
$ws.Range("E2").Value = 3930
$ws.Range("F2").Value = 721954
$ws.Range("G2").Value = 1652142
$ws.Range("H2").Value = -0.5630194014800181
$ws.Range("I2").Value = "TRENT"

# Row 3: JSL
$ws.Range("A3").Value = "JSL"
$ws.Range("B3").Value = 780.4
$ws.Range("C3").Value = 799.25
$ws.Range("D3").Value = 769.55
$ws.Range("E3").Value = 795
$ws.Range("F3").Value = 752362
$ws.Range("G3").Value = 1783303
$ws.Range("H3").Value = -0.5781075902412546
$ws.Range("I3").Value = "JSL"

# Row 4: MFSL
$ws.Range("A4").Value = "MFSL"
$ws.Range("B4").Value = 1631
$ws.Range("C4").Value = 1658.3
$ws.Range("D4").Value = 1627.6
$ws.Range("E4").Value = 1649.9
$ws.Range("F4").Value = 422995
$ws.Range("G4").Value = 888298
$ws.Range("H4").Value = -0.5238140804099525
$ws.Range("I4").Value = "MFSL"

# Row 5: GODREJPROP
$ws.Range("A5").Value = "GODREJPROP"
$ws.Range("B5").Value = 1875
$ws.Range("C5").Value = 1876.4
$ws.Range("D5").Value = 1850.1
$ws.Range("E5").Value = 1871.8
$ws.Range("F5").Value = 918048
$ws.Range("G5").Value = 2030980
$ws.Range("H5").Value = -0.5479778235137717
$ws.Range("I5").Value = "GODREJPROP"

# Row 6: NYKAA
$ws.Range("A6").Value = "NYKAA"
$ws.Range("B6").Value = 255.9
$ws.Range("C6").Value = 256.3
$ws.Range("D6").Value = 251.1
$ws.Range("E6").Value = 253.4
$ws.Range("F6").Value = 3678434
$ws.Range("G6").Value = 8428359
$ws.Range("H6").Value = -0.5635646274678143
$ws.Range("I6").Value = "NYKAA"

# Row 7: OBEROIRLTY
$ws.Range("A7").Value = "OBEROIRLTY"
$ws.Range("B7").Value = 1654
$ws.Range("C7").Value = 1666.9
$ws.Range("D7").Value = 1637.5
$ws.Range("E7").Value = 1646.7
$ws.Range("F7").Value = 212916
$ws.Range("G7").Value = 443321
$ws.Range("H7").Value = -0.5197249848304051
$ws.Range("I7").Value = "OBEROIRLTY"

# Row 8: ALKEM
$ws.Range("A8").Value = "ALKEM"
$ws.Range("B8").Value = 5888
$ws.Range("C8").Value = 5888.5
$ws.Range("D8").Value = 5778.5
$ws.Range("E8").Value = 5822.5
$ws.Range("F8").Value = 89063
$ws.Range("G8").Value = 183102
$ws.Range("H8").Value = -0.5135880547454424
$ws.Range("I8").Value = "ALKEM"

# Row 9: KFINTECH
$ws.Range("A9").Value = "KFINTECH"
$ws.Range("B9").Value = 1054.9
$ws.Range("C9").Value = 1074.2
$ws.Range("D9").Value = 1049.5
$ws.Range("E9").Value = 1070
$ws.Range("F9").Value = 472865
$ws.Range("G9").Value = 1065146
$ws.Range("H9").Value = -0.5560561650703284
$ws.Range("I9").Value = "KFINTECH"

# Row 10: PPLPHARMA
$ws.Range("A10").Value = "PPLPHARMA"
$ws.Range("B10").Value = 168.3
$ws.Range("C10").Value = 169.17
$ws.Range("D10").Value = 166.71
$ws.Range("E10").Value = 168.28
$ws.Range("F10").Value = 1416806
$ws.Range("G10").Value = 3404750
$ws.Range("H10").Value = -0.5838737058521184
$ws.Range("I10").Value = "PPLPHARMA"

# Row 11: INOXWIND
$ws.Range("A11").Value = "INOXWIND"
$ws.Range("B11").Value = 114.95
$ws.Range("C11").Value = 118.49
$ws.Range("D11").Value = 114.1
$ws.Range("E11").Value = 114.39
$ws.Range("F11").Value = 10797027
$ws.Range("G11").Value = 22393302
$ws.Range("H11").Value = -0.517845693323834
$ws.Range("I11").Value = "INOXWIND"

# Remove 2 now-unused trailing row(s) so the sheet shrinks back to A1:I11
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(12).Delete()


# ---- Pos_Change (sheet 2) ----
$ws = $wb.Worksheets.Item(2)

# Row 2: AXISBANK
$ws.Range("A2").Value = "AXISBANK"
$ws.Range("B2").Value = 1255
$ws.Range("C2").Value = 1308
$ws.Range("D2").Value = 1253.5
$ws.Range("E2").Value = 1299
$ws.Range("F2").Value = 7776704
$ws.Range("G2").Value = 5163166
$ws.Range("H2").Value = 0.5061890320783798
$ws.Range("I2").Value = "AXISBANK"

# Row 3: ONGC
$ws.Range("A3").Value = "ONGC"
$ws.Range("B3").Value = 243.99
$ws.Range("C3").Value = 251
$ws.Range("D3").Value = 243.58
$ws.Range("E3").Value = 247.98
$ws.Range("F3").Value = 33094088
$ws.Range("G3").Value = 22527192
$ws.Range("H3").Value = 0.4690729319481984
$ws.Range("I3").Value = "ONGC"

# Row 4: ULTRACEMCO
$ws.Range("A4").Value = "ULTRACEMCO"
$ws.Range("B4").Value = 12074
$ws.Range("C4").Value = 12323
$ws.Range("D4").Value = 11942
$ws.Range("E4").Value = 12242
$ws.Range("F4").Value = 408660
$ws.Range("G4").Value = 290750
$ws.Range("H4").Value = 0.4055374032674119
$ws.Range("I4").Value = "ULTRACEMCO"

# Row 5: JSWSTEEL
$ws.Range("A5").Value = "JSWSTEEL"
$ws.Range("B5").Value = 1179.9
$ws.Range("C5").Value = 1201.9
$ws.Range("D5").Value = 1170.3
$ws.Range("E5").Value = 1185.7
$ws.Range("F5").Value = 2281726
$ws.Range("G5").Value = 1570887
$ws.Range("H5").Value = 0.4525080416350762
$ws.Range("I5").Value = "JSWSTEEL"

# Row 6: BAJAJ-AUTO
$ws.Range("A6").Value = "BAJAJ-AUTO"
$ws.Range("B6").Value = 9576
$ws.Range("C6").Value = 9610
$ws.Range("D6").Value = 9494
$ws.Range("E6").Value = 9570
$ws.Range("F6").Value = 335206
$ws.Range("G6").Value = 229215
$ws.Range("H6").Value = 0.4624086556289946
$ws.Range("I6").Value = "BAJAJ-AUTO"

# Row 7: POWERGRID
$ws.Range("A7").Value = "POWERGRID"
$ws.Range("B7").Value = 257.1
$ws.Range("C7").Value = 259.7
$ws.Range("D7").Value = 256.75
$ws.Range("E7").Value = 258.2
$ws.Range("F7").Value = 15435341
$ws.Range("G7").Value = 9734559
$ws.Range("H7").Value = 0.5856230364416097
$ws.Range("I7").Value = "POWERGRID"

# Row 8: HDFCLIFE
$ws.Range("A8").Value = "HDFCLIFE"
$ws.Range("B8").Value = 742
$ws.Range("C8").Value = 753.75
$ws.Range("D8").Value = 741.85
$ws.Range("E8").Value = 743.9
$ws.Range("F8").Value = 2909648
$ws.Range("G8").Value = 1820910
$ws.Range("H8").Value = 0.5979087379387229
$ws.Range("I8").Value = "HDFCLIFE"

# Row 9: SUNPHARMA
$ws.Range("A9").Value = "SUNPHARMA"
$ws.Range("B9").Value = 1737.1
$ws.Range("C9").Value = 1737.1
$ws.Range("D9").Value = 1689.2
$ws.Range("E9").Value = 1706.3
$ws.Range("F9").Value = 2479677
$ws.Range("G9").Value = 1770082
$ws.Range("H9").Value = 0.4008825579831895
$ws.Range("I9").Value = "SUNPHARMA"

# Row 10: AMBUJACEM
$ws.Range("A10").Value = "AMBUJACEM"
$ws.Range("B10").Value = 537.65
$ws.Range("C10").Value = 551.35
$ws.Range("D10").Value = 533.5
$ws.Range("E10").Value = 549.5
$ws.Range("F10").Value = 1736723
$ws.Range("G10").Value = 1093928
$ws.Range("H10").Value = 0.5876026575789266
$ws.Range("I10").Value = "AMBUJACEM"

# Row 11: ZYDUSLIFE
$ws.Range("A11").Value = "ZYDUSLIFE"
$ws.Range("B11").Value = 900.7
$ws.Range("C11").Value = 902
$ws.Range("D11").Value = 881
$ws.Range("E11").Value = 885
$ws.Range("F11").Value = 728209
$ws.Range("G11").Value = 493965
$ws.Range("H11").Value = 0.4742117356492869
$ws.Range("I11").Value = "ZYDUSLIFE"

# Row 12: NATIONALUM
$ws.Range("A12").Value = "NATIONALUM"
$ws.Range("B12").Value = 359
$ws.Range("C12").Value = 374.3
$ws.Range("D12").Value = 356.3
$ws.Range("E12").Value = 373.75
$ws.Range("F12").Value = 27306024
$ws.Range("G12").Value = 17824247
$ws.Range("H12").Value = 0.5319594707142468
$ws.Range("I12").Value = "NATIONALUM"

# Row 13: INDUSTOWER
$ws.Range("A13").Value = "INDUSTOWER"
$ws.Range("B13").Value = 432.25
$ws.Range("C13").Value = 443.7
$ws.Range("D13").Value = 432.25
$ws.Range("E13").Value = 438.95
$ws.Range("F13").Value = 9332296
$ws.Range("G13").Value = 5850975
$ws.Range("H13").Value = 0.5949984404308684
$ws.Range("I13").Value = "INDUSTOWER"

# Row 14: RVNL
$ws.Range("A14").Value = "RVNL"
$ws.Range("B14").Value = 329.5
$ws.Range("C14").Value = 341.9
$ws.Range("D14").Value = 328
$ws.Range("E14").Value = 338.55
$ws.Range("F14").Value = 10947879
$ws.Range("G14").Value = 7141388
$ws.Range("H14").Value = 0.5330183712185922
$ws.Range("I14").Value = "RVNL"

# Row 15: PERSISTENT
$ws.Range("A15").Value = "PERSISTENT"
$ws.Range("B15").Value = 6325
$ws.Range("C15").Value = 6350
$ws.Range("D15").Value = 6230
$ws.Range("E15").Value = 6300
$ws.Range("F15").Value = 327713
$ws.Range("G15").Value = 228016
$ws.Range("H15").Value = 0.437236860571188
$ws.Range("I15").Value = "PERSISTENT"

# Row 16: LAURUSLABS
$ws.Range("A16").Value = "LAURUSLABS"
$ws.Range("B16").Value = 1063
$ws.Range("C16").Value = 1098.9
$ws.Range("D16").Value = 1057.8
$ws.Range("E16").Value = 1090.4
$ws.Range("F16").Value = 2383234
$ws.Range("G16").Value = 1574593
$ws.Range("H16").Value = 0.5135555664225613
$ws.Range("I16").Value = "LAURUSLABS"

# Row 17: CYIENT
$ws.Range("A17").Value = "CYIENT"
$ws.Range("B17").Value = 1175
$ws.Range("C17").Value = 1185.2
$ws.Range("D17").Value = 1161.1
$ws.Range("E17").Value = 1173.1
$ws.Range("F17").Value = 255742
$ws.Range("G17").Value = 177933
$ws.Range("H17").Value = 0.4372938128396643
$ws.Range("I17").Value = "CYIENT"


